$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F6").Value = "PR_A_Y1"
